# Refitting NCDEs to individual patients (for manuscript figure)
# - Updates the refit Prediction/Error values (columns D/E, plus F/G on the
#   last row of each Iterations block) to the latest per-patient fit results.
# - Adds a new "Label" column (H) recording the ground-truth class (0 =
#   Control, 1 = MDD) for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, matching the existing bold/centered/bordered
# header style used by B1:G1
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

# Row 2
$ws.Range("D2").Value = 0.6529379120031147
$ws.Range("E2").Value = 0.6529379120031147
$ws.Range("H2").Value = 0

# Row 3
$ws.Range("D3").Value = 0.5540253806261308
$ws.Range("E3").Value = 0.5540253806261308
$ws.Range("H3").Value = 0

# Row 4
$ws.Range("D4").Value = 0.5561251131271631
$ws.Range("E4").Value = 0.5561251131271631
$ws.Range("H4").Value = 0

# Row 5
$ws.Range("D5").Value = 0.5683137440666566
$ws.Range("E5").Value = 0.5683137440666566
$ws.Range("H5").Value = 0

# Row 6
$ws.Range("D6").Value = 0.5817466697773838
$ws.Range("E6").Value = 0.5817466697773838
$ws.Range("H6").Value = 0

# Row 7
$ws.Range("D7").Value = 0.5549218522105053
$ws.Range("E7").Value = 0.4450781477894947
$ws.Range("H7").Value = 1

# Row 8
$ws.Range("D8").Value = 0.5495809375530704
$ws.Range("E8").Value = 0.4504190624469296
$ws.Range("H8").Value = 1

# Row 9
$ws.Range("D9").Value = 0.5408531116506093
$ws.Range("E9").Value = 0.4591468883493907
$ws.Range("H9").Value = 1

# Row 10
$ws.Range("D10").Value = 0.5803070023896501
$ws.Range("E10").Value = 0.4196929976103499
$ws.Range("H10").Value = 1

# Row 11
$ws.Range("D11").Value = 0.5501619079018494
$ws.Range("E11").Value = 0.4498380920981506
$ws.Range("F11").Value = 0.733355700969696
$ws.Range("G11").Value = 0.5
$ws.Range("H11").Value = 1

# Row 12
$ws.Range("D12").Value = 0.6669846687040633
$ws.Range("E12").Value = 0.6669846687040633
$ws.Range("H12").Value = 0

# Row 13
$ws.Range("D13").Value = 0.522628377530662
$ws.Range("E13").Value = 0.522628377530662
$ws.Range("H13").Value = 0

# Row 14
$ws.Range("D14").Value = 0.5337799252556555
$ws.Range("E14").Value = 0.5337799252556555
$ws.Range("H14").Value = 0

# Row 15
$ws.Range("D15").Value = 0.5383336366274961
$ws.Range("E15").Value = 0.5383336366274961
$ws.Range("H15").Value = 0

# Row 16
$ws.Range("D16").Value = 0.6116973178939507
$ws.Range("E16").Value = 0.6116973178939507
$ws.Range("H16").Value = 0

# Row 17
$ws.Range("D17").Value = 0.5068862471055573
$ws.Range("E17").Value = 0.4931137528944427
$ws.Range("H17").Value = 1

# Row 18
$ws.Range("D18").Value = 0.5117428467115197
$ws.Range("E18").Value = 0.4882571532884803
$ws.Range("H18").Value = 1

# Row 19
$ws.Range("D19").Value = 0.4804063889058778
$ws.Range("E19").Value = 0.5195936110941222
$ws.Range("H19").Value = 1

# Row 20
$ws.Range("D20").Value = 0.6290948522813512
$ws.Range("E20").Value = 0.3709051477186488
$ws.Range("H20").Value = 1

# Row 21
$ws.Range("D21").Value = 0.4991523018633796
$ws.Range("E21").Value = 0.5008476981366203
$ws.Range("F21").Value = 0.7561848759651184
$ws.Range("G21").Value = 0.3
$ws.Range("H21").Value = 1
